# Insert a new weekly data row for "Terminal Hortofrutícola Agro Chillán" (Mango)
# at row 80, shifting the existing rows 80-98 down to 81-99 (matches the target
# diff: dimension grows from A1:T98 to A1:T99, and the new record lands at row 80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 80 downward by one row.
$ws.Rows.Item(80).Insert(-4121)   # -4121 = xlShiftDown

# Populate the newly-opened row 80 with the new observation.
$ws.Cells.Item(80, 1).Value2 = 7
$ws.Cells.Item(80, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(80, 3).Value2 = "Ñuble"
$ws.Cells.Item(80, 4).Value2 = 44985
$ws.Cells.Item(80, 5).Value2 = 16
$ws.Cells.Item(80, 6).Value2 = "Fruta"
$ws.Cells.Item(80, 7).Value2 = 100108
$ws.Cells.Item(80, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(80, 9).Value2 = 100108002
$ws.Cells.Item(80, 10).Value2 = "Mango"
$ws.Cells.Item(80, 11).Value2 = "Sin especificar"
$ws.Cells.Item(80, 12).Value2 = "Primera"
$ws.Cells.Item(80, 13).Value2 = 60
$ws.Cells.Item(80, 14).Value2 = 8000
$ws.Cells.Item(80, 15).Value2 = 8000
$ws.Cells.Item(80, 16).Value2 = 8000
$ws.Cells.Item(80, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(80, 18).Value2 = "Perú"
$ws.Cells.Item(80, 19).Value2 = 2000
$ws.Cells.Item(80, 20).Value2 = 4
